$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.339.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.84%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.766.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.69%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'614.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.18%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'178.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.34%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.763.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.69%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.98%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D13").Value = "'39.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.76%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -3.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.395.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.66%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.766.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'69.434.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.74%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.50%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -3.29%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'508.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.30%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'16.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.16%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.40%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.30%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'86.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "'  -3.09%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'RenderToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'10.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.31%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'PEPE"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.0000135"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.89%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.64%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +3.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'8.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.94%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'30.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.35%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.09%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.43%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.47%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.43%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'456.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +8.60%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'49.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.37%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +5.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.62%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.60%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.958.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.70%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.54%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'27.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.33%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'139.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.24%  "
$ws.Range("E51").Style = "Normal"
